$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "[$-409]d\-mmm\-yyyy;@"

# New wallet-ledger transactions for 1-2 Mar 2020, entered into the rows
# that were previously blank placeholders (80-83).
$ws.Range("A80").Value = 43891
$ws.Range("A80").NumberFormat = $dateFormat
$ws.Range("B80").Value = 44720
$ws.Range("D80").Value = "Ordered Amount"
$ws.Range("E80").Formula = '=IF(A80="","",SUM(E79-B80+C80))'

$ws.Range("A81").Value = 43891
$ws.Range("A81").NumberFormat = $dateFormat
$ws.Range("C81").Value = 42236
$ws.Range("D81").Value = "Manual Added"
$ws.Range("E81").Formula = '=IF(A81="","",SUM(E80-B81+C81))'

$ws.Range("A82").Value = 43892
$ws.Range("A82").NumberFormat = $dateFormat
$ws.Range("C82").Value = 40018
$ws.Range("D82").Value = "Manual Added"
$ws.Range("E82").Formula = '=IF(A82="","",SUM(E81-B82+C82))'

$ws.Range("A83").Value = 43892
$ws.Range("A83").NumberFormat = $dateFormat
$ws.Range("B83").Value = 29120
$ws.Range("D83").Value = "Ordered Amount"
$ws.Range("E83").Formula = '=IF(A83="","",SUM(E82-B83+C83))'

# Extend the running-balance fill-down formula by one more row so the
# sheet keeps the same buffer of blank template rows below the data.
$ws.Range("E151").Formula = '=IF(A151="","",SUM(E150-B151+C151))'

$ws.Range("E83").Select()
